$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 68.06
$ws.Range("B3").Value = 73.76000000000001
$ws.Range("B4").Value = 61.56
$ws.Range("B5").Value = 72.7
$ws.Range("B6").Value = 53.66
$ws.Range("B7").Value = 63.9
$ws.Range("B8").Value = 57.45
$ws.Range("B9").Value = 62.56
$ws.Range("B10").Value = 73.95999999999999
$ws.Range("B11").Value = 51.56
$ws.Range("B12").Value = 84.26000000000001
$ws.Range("B13").Value = 43.51
$ws.Range("B14").Value = 50.26
$ws.Range("B15").Value = 59.25
$ws.Range("B16").Value = 60.16
$ws.Range("B17").Value = 64.2
$ws.Range("B19").Value = 58.96
$ws.Range("B20").Value = 68.06
$ws.Range("B21").Value = 51.26
